# "new cue style and new data dir name"
#
# Appends a new subject (S4 / Nick, session dir "2013-03-04-nick") with 15
# new runs recorded on 2013-03-04 to the watchErpDataset2 log sheet, and
# moves the sheet's row/selection bookkeeping down to match.
#
# Columns: A subjectTag | B subjectName | C date | D sessionDirectory
#          E fileName | F condition | G run

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$lastRow = 46
$firstNewRow = 47
$lastNewRow = 61

# ---------------------------------------------------------------------
# 1) Clone the formatting (cell styles / number formats) of the last
#    existing data row down across the whole new block in one shot, so
#    the new rows pick up the same (centered / date) styles already
#    registered in the stylesheet instead of minting new ones.
# ---------------------------------------------------------------------
$ws.Range("A$lastRow`:G$lastRow").Copy() | Out-Null
$ws.Range("A$firstNewRow`:G$lastNewRow").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) New row data, 2013-03-04, subject S4 / Nick.
# ---------------------------------------------------------------------
$fileNames = @(
    "2013-03-04-14-51-57",
    "2013-03-04-15-00-40",
    "2013-03-04-15-06-10",
    "2013-03-04-15-12-15",
    "2013-03-04-15-21-41",
    "2013-03-04-15-34-23",
    "2013-03-04-15-40-19",
    "2013-03-04-15-46-20",
    "2013-03-04-15-52-05",
    "2013-03-04-15-58-31",
    "2013-03-04-16-16-00",
    "2013-03-04-16-21-35",
    "2013-03-04-16-29-47",
    "2013-03-04-16-36-52",
    "2013-03-04-16-44-05"
)

$conditions = @(
    "hybrid-10Hz",
    "oddball.bdf",
    "hybrid-12Hz",
    "hybrid-10Hz",
    "hybrid-15Hz",
    "hybrid-10Hz",
    "hybrid-12Hz",
    "oddball.bdf",
    "hybrid-8-57Hz",
    "oddball.bdf",
    "hybrid-15Hz",
    "hybrid-15Hz",
    "hybrid-8-57Hz",
    "hybrid-12Hz",
    "hybrid-8-57Hz"
)

$runs = @(1, 1, 1, 2, 1, 3, 2, 2, 1, 3, 2, 3, 2, 3, 3)

# Column A (subjectTag) - identical for every new row.
for ($r = $firstNewRow; $r -le $lastNewRow; $r++) {
    $ws.Cells.Item($r, 1).Value = "S4"
}

# Column B (subjectName) - identical for every new row.
for ($r = $firstNewRow; $r -le $lastNewRow; $r++) {
    $ws.Cells.Item($r, 2).Value = "Nick"
}

# Column C (date) - 2013-03-04 for every new row.
for ($r = $firstNewRow; $r -le $lastNewRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 41337
}

# Column E (fileName) - unique timestamp per row.
for ($r = $firstNewRow; $r -le $lastNewRow; $r++) {
    $ws.Cells.Item($r, 5).Value = $fileNames[$r - $firstNewRow]
}

# Column D (sessionDirectory) - new data-dir name, identical for every new row.
for ($r = $firstNewRow; $r -le $lastNewRow; $r++) {
    $ws.Cells.Item($r, 4).Value = "2013-03-04-nick"
}

# Column F (condition / cue style) per row.
for ($r = $firstNewRow; $r -le $lastNewRow; $r++) {
    $ws.Cells.Item($r, 6).Value = $conditions[$r - $firstNewRow]
}

# Column G (run) per row.
for ($r = $firstNewRow; $r -le $lastNewRow; $r++) {
    $ws.Cells.Item($r, 7).Value = $runs[$r - $firstNewRow]
}

# ---------------------------------------------------------------------
# 3) Move the sheet's selection/scroll bookkeeping down to the new block,
#    mirroring how the author's view state shifted after appending rows.
# ---------------------------------------------------------------------
$ws.Range("A$($firstNewRow + 1):A$lastNewRow").Select()
try { $excel.ActiveWindow.ScrollRow = 22 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 1 } catch {}
